$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "001762-27-2<br><span style='font-size:smaller'>diethyldimethylplumbane</span>"
$ws.Cells.Item(2, 2).Value = -1.864285362608939
$ws.Cells.Item(2, 3).Value = 0.6933238970689278
$ws.Cells.Item(2, 4).Value = 0.007168580467208378
$ws.Cells.Item(2, 5).Value = "**"
$ws.Cells.Item(2, 6).Value = 0.1550069440709066
$ws.Cells.Item(2, 7).Value = 0.03982739842619181
$ws.Cells.Item(2, 8).Value = 0.6032820033356766

# Row 3
$ws.Cells.Item(3, 1).Value = "Intercept"
$ws.Cells.Item(3, 2).Value = 0.9791239768184468
$ws.Cells.Item(3, 3).Value = 0.4655359435894203
$ws.Cells.Item(3, 4).Value = 0.03544664234137789
$ws.Cells.Item(3, 5).Value = "*"
$ws.Cells.Item(3, 6).Value = 2.662123138570772
$ws.Cells.Item(3, 7).Value = 1.068946439460581
$ws.Cells.Item(3, 8).Value = 6.629798597289993

# Row 4
$ws.Cells.Item(4, 1).Value = "1000401-12-0<br><span style='font-size:smaller'>2,5-cyclohexadien-1-one, 2,6-bis(1,1-dimethylethyl)-4-hydroxy-4-methyl-</span>"
$ws.Cells.Item(4, 2).Value = -0.5446317316507132
$ws.Cells.Item(4, 3).Value = 0.3947848536714565
$ws.Cells.Item(4, 4).Value = 0.1677203352096722
$ws.Cells.Item(4, 5).Value = ""
$ws.Cells.Item(4, 6).Value = 0.5800553600406002
$ws.Cells.Item(4, 7).Value = 0.2675603729481738
$ws.Cells.Item(4, 8).Value = 1.25752635565732

# Row 5
$ws.Cells.Item(5, 1).Value = "019549-87-2<br><span style='font-size:smaller'>2,4-Dimethyl-1-heptene</span>"
$ws.Cells.Item(5, 2).Value = -0.382333806875028
$ws.Cells.Item(5, 3).Value = 0.4609742337128525
$ws.Cells.Item(5, 4).Value = 0.4068759203773512
$ws.Cells.Item(5, 5).Value = ""
$ws.Cells.Item(5, 6).Value = 0.6822672696844253
$ws.Cells.Item(5, 7).Value = 0.2764173811689787
$ws.Cells.Item(5, 8).Value = 1.684006357755336

# Row 6
$ws.Cells.Item(6, 1).Value = "000109-52-4<br><span style='font-size:smaller'>Pentanoic acid</span>"
$ws.Cells.Item(6, 2).Value = -0.2684651426305072
$ws.Cells.Item(6, 3).Value = 0.5365811253559523
$ws.Cells.Item(6, 4).Value = 0.6168460028724494
$ws.Cells.Item(6, 5).Value = ""
$ws.Cells.Item(6, 6).Value = 0.7645520726194696
$ws.Cells.Item(6, 7).Value = 0.2670914557513763
$ws.Cells.Item(6, 8).Value = 2.188538267172609

# Row 7
$ws.Cells.Item(7, 1).Value = "063521-76-6<br><span style='font-size:smaller'>Tridecanedial</span>"
$ws.Cells.Item(7, 2).Value = 0.06864347057661022
$ws.Cells.Item(7, 3).Value = 0.3921975995195186
$ws.Cells.Item(7, 4).Value = 0.861061821359727
$ws.Cells.Item(7, 5).Value = ""
$ws.Cells.Item(7, 6).Value = 1.071054278703266
$ws.Cells.Item(7, 7).Value = 0.4965536105032165
$ws.Cells.Item(7, 8).Value = 2.310238499254941
